$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 (description) texts in-place before shifting columns
$ws.Range("C1").Value = "Cannot be cleared"
$ws.Range("E1").Value = "To clear, set it to NULL"
$ws.Range("F1").Value = "Comma-delimited list of tags - This will replace the current list of tags assigned to the application. To clear the tags, set it to NULL"
$ws.Range("G1").Value = "To clear the Business Unit, set it to NULL"
$ws.Range("H1").Value = "To clear the Business Owner, set it to NULL"
$ws.Range("I1").Value = "To clear the Owner Email, set it to NULL"
$ws.Range("J1").Value = "Comma-delimited list of teams - This will replace the current list of teams assigned to the application."
$ws.Range("M1").Value = "Custom fields can be renamed, in that case, you need to put their real name here. Keeping the original names will throw an error"

# Remove the "Dynamic Scan Approval" (K) and "Archer Application Name" (L) columns entirely,
# shifting everything to the right of them (including the custom fields note) to the left
$ws.Columns("K:L").Delete()

# Update the selected/active cell to match the saved view state
$ws.Range("J2").Select() | Out-Null
